$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G7").Value = '[[2,8,"","customer"]]'
$ws.Range("G8").Value = '[[2,9,"","peer"]]'
$ws.Range("G9").Value = '[[1,6,"","provider"]]'
$ws.Range("G10").Value = '[[1,7,"","peer"]]'
$ws.Range("E11").Value = '[[8,"",500],[9,""],[11,""],[12,"",10]]'
